$d = $word.ActiveDocument

# 1. Modify first paragraph: split "This is a Microsoft word document." into
#    the original text (with two trailing spaces) plus three new red runs.
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$end = $d.Paragraphs(1).Range.End
$ins = $d.Range($end - 1, $end - 1)
$ins.InsertAfter("(This is a change – Ve")
$ins.Font.Color = 255

$end2 = $d.Paragraphs(1).Range.End
$ins2 = $d.Range($end2 - 1, $end2 - 1)
$ins2.InsertAfter("rsion for main branch")
$ins2.Font.Color = 255

$end3 = $d.Paragraphs(1).Range.End
$ins3 = $d.Range($end3 - 1, $end3 - 1)
$ins3.InsertAfter(")")
$ins3.Font.Color = 255

# 2. Remove the last paragraph ("ank God almighty, we are free at last.")
$lastParaCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaCount)
$lastPara.Range.Delete()
